# Generate Report for Handback
# -----------------------------------------------------------------------
# This script reproduces a "handback" update on the localization-status
# workbook:
#   1. The Status text "In Translation" is renamed everywhere to
#      "Handed back: in sync with en-US" (Overview + zh-cn + de-de sheets).
#   2. On the zh-cn and de-de detail sheets, rows 2-4 get their
#      "Latest Target File" (J), "Latest Handback File" (K) and
#      "Latest Handback DateTime" (L) columns populated: J gets a
#      hyperlink (same target as column A's hyperlink) showing the
#      source file name, K gets the generated handback xliff file name,
#      and L gets the handback timestamp.
#   3. A couple of columns that now hold much longer text are widened.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/b787781a66dcf9b8df6172274a7d074c484ce7c7/e2e/"

# ---------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
#    (this is a single shared string, used by every Status cell across
#    all three sheets, so a global replace covers Overview + both
#    per-language sheets in one shot).
# ---------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("In Translation", "Handed back: in sync with en-US")
}

# ---------------------------------------------------------------------
# 2. Per-language handback details for the three source files.
# ---------------------------------------------------------------------
$files = @(
    @{ Row = 2; Name = "14f64019-da68-4340-9233-f9329e56b714.yml"; Hash = "14f64019-da68-4340-9233-f9329e56b714.76b4bf4e2c2d897cd93a205b8e037d8cf880db72" },
    @{ Row = 3; Name = "729e5e2f-2429-4344-bdec-9d2d67b82b67.yml"; Hash = "729e5e2f-2429-4344-bdec-9d2d67b82b67.3da2f0f866f23e9a2f9fb1234b638ea118724000" },
    @{ Row = 4; Name = "e557a9b9-431d-48d4-8775-cfbe719ab9a2.md";  Hash = "e557a9b9-431d-48d4-8775-cfbe719ab9a2.3010c438b5d95f2aeb8f23e8bfab54d549f7abb4" }
)

$languages = @(
    @{ Sheet = "zh-cn"; Suffix = "zh-cn"; HandbackTime = "2017-05-22 02:09:37" },
    @{ Sheet = "de-de"; Suffix = "de-de"; HandbackTime = "2017-05-22 02:10:06" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    foreach ($f in $files) {
        $row = $f.Row
        $fileUrl = $urlBase + $f.Name

        # J: Latest Target File -- a hyperlink to the source file, same
        # target/display text as the column-A hyperlink on this row.
        $ws.Hyperlinks.Add($ws.Range("J" + $row), $fileUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $f.Name)

        # K: Latest Handback File -- generated xliff name for this language.
        $ws.Range("K" + $row).Value = $f.Hash + "." + $lang.Suffix + ".xlf"

        # L: Latest Handback DateTime.
        $ws.Range("L" + $row).Value = $lang.HandbackTime
    }
}

# ---------------------------------------------------------------------
# 3. Widen columns that now contain longer text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = 29.1666666666667
$overview.Range("F1").ColumnWidth = 29.1666666666667

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)
    $ws.Range("C1").ColumnWidth = 29.1666666666667
    $ws.Range("J1").ColumnWidth = 39.1666666666667
    $ws.Range("K1").ColumnWidth = 39.1666666666667
}
